$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2944.3333
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2944.3333
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2944.3333
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -3912.3333
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 2500
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2500
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2252
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 2500
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2500
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1642
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 2833.3333
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -4996
$ws.Range("H89").Value = 2833.3333
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -24982
$ws.Range("H98").Value = 432044.12
$ws.Range("I98").Value = 590096.06
$ws.Range("J98").Value = 3046
$ws.Range("K98").Value = 590096.06
$ws.Range("L98").Value = 3046
$ws.Range("M98").Value = -588598.06
$ws.Range("N98").Value = -6042
$ws.Range("H106").Value = 4275914.5
$ws.Range("I106").Value = 4833255.5
$ws.Range("K106").Value = 4833255.5
$ws.Range("M106").Value = -4832624.5
$ws.Range("H122").Value = 432044.12
$ws.Range("I122").Value = 590096.06
$ws.Range("J122").Value = 3046
$ws.Range("K122").Value = 1770288.18
$ws.Range("L122").Value = 9138
$ws.Range("M122").Value = -1767838.18
$ws.Range("N122").Value = -14038
$ws.Range("H138").Value = 2544.6
$ws.Range("I138").Value = 898.7778
$ws.Range("J138").Value = 3153.3289
$ws.Range("K138").Value = 2696.3334
$ws.Range("L138").Value = 9459.9867
$ws.Range("M138").Value = 2443.6666
$ws.Range("N138").Value = -19739.9867
$ws.Range("H141").Value = 2391.3572
$ws.Range("I141").Value = 1495.186
$ws.Range("K141").Value = 4485.558
$ws.Range("M141").Value = 694.442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18409.758
$ws.Range("I32").Value = 2248.1833
$ws.Range("K32").Value = 2248.1833
$ws.Range("M32").Value = -1961.1833
$ws.Range("H45").Value = 890
$ws.Range("I45").Value = 862.5
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 862.5
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -485.5
$ws.Range("N45").Value = -1754
$ws.Range("H61").Value = 2231.05
$ws.Range("I61").Value = 1495.1613
$ws.Range("J61").Value = 4765.778
$ws.Range("K61").Value = 1495.1613
$ws.Range("L61").Value = 4765.778
$ws.Range("M61").Value = -1283.1613
$ws.Range("N61").Value = -5189.778
$ws.Range("H102").Value = 1758.7142
$ws.Range("I102").Value = 1860
$ws.Range("J102").Value = 1505.5
$ws.Range("K102").Value = 1860
$ws.Range("L102").Value = 1505.5
$ws.Range("M102").Value = -238
$ws.Range("N102").Value = -4749.5
$ws.Range("H132").Value = 1718.4445
$ws.Range("I132").Value = 1332.4166
$ws.Range("J132").Value = 4806.6665
$ws.Range("K132").Value = 3997.2498
$ws.Range("L132").Value = 14419.9995
$ws.Range("M132").Value = -1467.2498
$ws.Range("N132").Value = -19479.9995
$ws.Range("H136").Value = 2231.05
$ws.Range("I136").Value = 1495.1613
$ws.Range("J136").Value = 4765.778
$ws.Range("K136").Value = 4485.4839
$ws.Range("L136").Value = 14297.334
$ws.Range("M136").Value = -1935.4839
$ws.Range("N136").Value = -19397.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1880.6923
$ws.Range("I20").Value = 1919.9259
$ws.Range("J20").Value = 1792.4166
$ws.Range("K20").Value = 1919.9259
$ws.Range("L20").Value = 1792.4166
$ws.Range("M20").Value = -1672.9259
$ws.Range("N20").Value = -2286.4166
$ws.Range("H59").Value = 41775
$ws.Range("I59").Value = 20000
$ws.Range("K59").Value = 20000
$ws.Range("M59").Value = -19153
$ws.Range("H86").Value = 5904.391
$ws.Range("I86").Value = 1706.2142
$ws.Range("J86").Value = 12434.889
$ws.Range("K86").Value = 1706.2142
$ws.Range("L86").Value = 12434.889
$ws.Range("M86").Value = -583.2141999999999
$ws.Range("N86").Value = -14680.889
$ws.Range("H89").Value = 5904.391
$ws.Range("I89").Value = 1706.2142
$ws.Range("J89").Value = 12434.889
$ws.Range("K89").Value = 8531.071
$ws.Range("L89").Value = 62174.44499999999
$ws.Range("M89").Value = -2915.071
$ws.Range("N89").Value = -73406.44499999999
$ws.Range("H94").Value = 1136.5714
$ws.Range("I94").Value = 1016.3571
$ws.Range("J94").Value = 1377
$ws.Range("K94").Value = 1016.3571
$ws.Range("L94").Value = 1377
$ws.Range("M94").Value = -565.3571
$ws.Range("N94").Value = -2279
$ws.Range("H105").Value = 2976.8684
$ws.Range("I105").Value = 2734.577
$ws.Range("J105").Value = 3501.8333
$ws.Range("K105").Value = 2734.577
$ws.Range("L105").Value = 3501.8333
$ws.Range("M105").Value = -987.5770000000002
$ws.Range("N105").Value = -6995.8333
$ws.Range("H132").Value = 30391.111
$ws.Range("J132").Value = 30391.111
$ws.Range("L132").Value = 30391.111
$ws.Range("N132").Value = -40511.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2598
$ws.Range("I132").Value = 1981.1034
$ws.Range("J132").Value = 4088.8333
$ws.Range("K132").Value = 5943.3102
$ws.Range("L132").Value = 12266.4999
$ws.Range("M132").Value = -3413.3102
$ws.Range("N132").Value = -17326.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1092.2142
$ws.Range("I23").Value = 3401
$ws.Range("J23").Value = 462.54544
$ws.Range("K23").Value = 10203
$ws.Range("L23").Value = 1387.63632
$ws.Range("M23").Value = -9968
$ws.Range("N23").Value = -1857.63632
$ws.Range("H117").Value = 1540
$ws.Range("I117").Value = 80
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 240
$ws.Range("L117").Value = 9000
$ws.Range("M117").Value = 3202
$ws.Range("N117").Value = -15884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3705339.2
$ws.Range("I122").Value = 5557059
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 16671177
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -16668727
$ws.Range("N122").Value = -10600
$ws.Range("H132").Value = 2575.0962
$ws.Range("I132").Value = 2420.1191
$ws.Range("J132").Value = 3226
$ws.Range("K132").Value = 7260.3573
$ws.Range("L132").Value = 9678
$ws.Range("M132").Value = -4730.3573
$ws.Range("N132").Value = -14738
$ws.Range("H135").Value = 166702500
$ws.Range("J135").Value = 166702500
$ws.Range("L135").Value = 166702500
$ws.Range("N135").Value = -166712640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2706
$ws.Range("I7").Value = 1936.4615
$ws.Range("J7").Value = 3331.25
$ws.Range("K7").Value = 1936.4615
$ws.Range("L7").Value = 3331.25
$ws.Range("M7").Value = -1824.4615
$ws.Range("N7").Value = -3555.25
$ws.Range("H16").Value = 581.8823
$ws.Range("I16").Value = 647.6667
$ws.Range("J16").Value = 424
$ws.Range("K16").Value = 647.6667
$ws.Range("L16").Value = 424
$ws.Range("M16").Value = -477.6667
$ws.Range("N16").Value = -764
$ws.Range("H93").Value = 2165
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2165
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 2165
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -4661
$ws.Range("H100").Value = 2434.4783
$ws.Range("I100").Value = 1856.2858
$ws.Range("J100").Value = 2687.4375
$ws.Range("K100").Value = 1856.2858
$ws.Range("L100").Value = 2687.4375
$ws.Range("M100").Value = -1315.2858
$ws.Range("N100").Value = -3769.4375
$ws.Range("H115").Value = 24814.4
$ws.Range("J115").Value = 24814.4
$ws.Range("L115").Value = 24814.4
$ws.Range("N115").Value = -27164.4
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H126").Value = 2706
$ws.Range("I126").Value = 1936.4615
$ws.Range("J126").Value = 3331.25
$ws.Range("K126").Value = 5809.3845
$ws.Range("L126").Value = 9993.75
$ws.Range("M126").Value = -3339.3845
$ws.Range("N126").Value = -14933.75
$ws.Range("H136").Value = 5453.143
$ws.Range("I136").Value = 3477.2727
$ws.Range("J136").Value = 12698
$ws.Range("K136").Value = 10431.8181
$ws.Range("L136").Value = 38094
$ws.Range("M136").Value = -7881.8181
$ws.Range("N136").Value = -43194

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 764
$ws.Range("I100").Value = 292
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 584
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -43
$ws.Range("N100").Value = -3082
$ws.Range("H116").Value = 72180
$ws.Range("J116").Value = 72180
$ws.Range("L116").Value = 72180
$ws.Range("N116").Value = -81358
$ws.Range("H132").Value = 9806765
$ws.Range("I132").Value = 12502786
$ws.Range("K132").Value = 37508358
$ws.Range("M132").Value = -37505828
$ws.Range("H136").Value = 7961180.5
$ws.Range("I136").Value = 9036889
$ws.Range("J136").Value = 940
$ws.Range("K136").Value = 27110667
$ws.Range("L136").Value = 2820
$ws.Range("M136").Value = -27108117
$ws.Range("N136").Value = -7920
